# Actualización automática 2025-09-10 17:15:09
#
# Updates the "septiembre" sales figures for client TULCAN NARVAEZ EDITH
# MARITZA (advisor HIDALGO HIDALGO PEDRO GUSTAVO): two product groups
# (INODOROS, LAVABOS) that were previously 0 now carry negative
# adjustments, and the dependent monthly / compliance totals are
# refreshed to stay consistent.

$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO": per-group sales for septiembre (row 22) ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Cells.Item(22, 8).Value = -67.65000000000001   # H22 INODOROS
$ws1.Cells.Item(22, 9).Value = -57.6                # I22 LAVABOS

# --- Sheet "VENTA MENSUAL": septiembre column totals ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Cells.Item(22, 6).Value = 7229.68     # F22 client row total for septiembre
$ws2.Cells.Item(23, 6).Value = 24077.51    # F23 grand total for septiembre

# --- Sheet "CUMPLIMIENTO MENSUAL": VENTA / POR CUMPLIR / CUMPLIMIENTO ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 6 = INODOROS
$ws3.Cells.Item(6, 4).Value = 1735.95                 # D6 VENTA
$ws3.Cells.Item(6, 5).Value = 1171.63368146026        # E6 POR CUMPLIR
$ws3.Cells.Item(6, 6).Value = 0.5970421457064181      # F6 CUMPLIMIENTO

# Row 7 = LAVABOS
$ws3.Cells.Item(7, 4).Value = 70.2                    # D7 VENTA
$ws3.Cells.Item(7, 5).Value = 816.511016287574        # E7 POR CUMPLIR
$ws3.Cells.Item(7, 6).Value = 0.07916897242791564     # F7 CUMPLIMIENTO

# Row 15 = TOTAL
$ws3.Cells.Item(15, 4).Value = 24077.51               # D15 VENTA
$ws3.Cells.Item(15, 5).Value = 31347.23316613378      # E15 POR CUMPLIR
$ws3.Cells.Item(15, 6).Value = 0.4344180707852535     # F15 CUMPLIMIENTO

# Column E on this sheet narrows by one character (23 -> 22) as a side
# effect of the refreshed "POR CUMPLIR" figures. ColumnWidth is specified
# in characters and the stored OOXML width is offset by 5/6 of a
# character from the COM property value, so back the input off by that
# amount to land exactly on the target stored width of 22.
$ws3.Columns.Item(5).ColumnWidth = 22 - (5/6)
